# Update the "Age_(days)" column (H) with the recalculated ages for the
# period that was chosen, and drop the old explicit centered style on
# those cells (they revert to the default/general cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray empty column that was sitting between the data (H) and the
# leftover bestFit column formatting (previously at J), so the latter slides
# back to I.
$ws.Columns("I:I").Delete()

$ageValues = @(12, 14, 16, 18, 1, 14, 16, 18, 12, 14, 16, 18, 12, 14, 16, 12, 14, 16, 18, 12, 14, 16, 18, 12, 14, 16, 18, 100)

for ($i = 0; $i -lt $ageValues.Length; $i++) {
    $row = 2 + $i
    $ws.Range("H$row").Value = $ageValues[$i]
}

# The previously applied centered number style is cleared from the data
# cells (header H1 keeps its style).
$ws.Range("H2:H29").Style = "Normal"

# Update the remembered selection left over from editing.
$ws.Range("L11").Select()
